$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1732954545454546
$ws.Range("C2").Value = 0.5852272727272727
$ws.Range("J2").Value = 0.008522727272727272
$ws.Range("P2").Value = 0.1392045454545454
$ws.Range("S2").Value = 0.09375
$ws.Range("C3").Value = 0.01428571428571429
$ws.Range("J3").Value = 0.02380952380952381
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2476190476190476
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.7073170731707317
$ws.Range("S4").Value = 0.2195121951219512
$ws.Range("B6").Value = 0.05985915492957746
$ws.Range("D6").Value = 0.01408450704225352
$ws.Range("E6").Value = 0.00352112676056338
$ws.Range("F6").Value = 0.05985915492957746
$ws.Range("J6").Value = 0.2922535211267606
$ws.Range("O6").Value = 0.0176056338028169
$ws.Range("Q6").Value = 0.1056338028169014
$ws.Range("R6").Value = 0.06690140845070422
$ws.Range("S6").Value = 0.3802816901408451
$ws.Range("B7").Value = 0.09230769230769231
$ws.Range("D7").Value = 0.01025641025641026
$ws.Range("F7").Value = 0.05128205128205128
$ws.Range("J7").Value = 0.1282051282051282
$ws.Range("O7").Value = 0.02051282051282051
$ws.Range("Q7").Value = 0.1435897435897436
$ws.Range("R7").Value = 0.08717948717948718
$ws.Range("S7").Value = 0.4666666666666667
$ws.Range("B8").Value = 0.1140350877192982
$ws.Range("D8").Value = 0.01096491228070175
$ws.Range("F8").Value = 0.06798245614035088
$ws.Range("J8").Value = 0.08991228070175439
$ws.Range("O8").Value = 0.02850877192982456
$ws.Range("Q8").Value = 0.1578947368421053
$ws.Range("R8").Value = 0.07894736842105263
$ws.Range("S8").Value = 0.4517543859649123
$ws.Range("B9").Value = 0.1063829787234043
$ws.Range("D9").Value = 0.00425531914893617
$ws.Range("E9").Value = 0.00425531914893617
$ws.Range("F9").Value = 0.09361702127659574
$ws.Range("J9").Value = 0.1106382978723404
$ws.Range("O9").Value = 0.008510638297872341
$ws.Range("Q9").Value = 0.148936170212766
$ws.Range("R9").Value = 0.07659574468085106
$ws.Range("S9").Value = 0.4468085106382979
$ws.Range("B10").Value = 0.1348228043143297
$ws.Range("D10").Value = 0.02311248073959938
$ws.Range("E10").Value = 0.002311248073959939
$ws.Range("F10").Value = 0.08320493066255778
$ws.Range("J10").Value = 0.09707241910631741
$ws.Range("O10").Value = 0.01848998459167951
$ws.Range("Q10").Value = 0.1687211093990755
$ws.Range("R10").Value = 0.08166409861325115
$ws.Range("S10").Value = 0.3906009244992296
$ws.Range("G11").Value = 0.1697530864197531
$ws.Range("J11").Value = 0.08641975308641975
$ws.Range("K11").Value = 0.2345679012345679
$ws.Range("L11").Value = 0.5
$ws.Range("S11").Value = 0.009259259259259259
$ws.Range("G12").Value = 0.7160493827160493
$ws.Range("J12").Value = 0.2222222222222222
$ws.Range("K12").Value = 0.01234567901234568
$ws.Range("L12").Value = 0.01851851851851852
$ws.Range("S12").Value = 0.0308641975308642
$ws.Range("G13").Value = 0.62
$ws.Range("J13").Value = 0.34
$ws.Range("S13").Value = 0.04
$ws.Range("F15").Value = 0.02459016393442623
$ws.Range("H15").Value = 0.1311475409836066
$ws.Range("I15").Value = 0.09426229508196721
$ws.Range("J15").Value = 0.3360655737704918
$ws.Range("K15").Value = 0.04918032786885246
$ws.Range("M15").Value = 0.02049180327868852
$ws.Range("O15").Value = 0.09836065573770492
$ws.Range("S15").Value = 0.2459016393442623
$ws.Range("F16").Value = 0.013215859030837
$ws.Range("H16").Value = 0.1762114537444934
$ws.Range("I16").Value = 0.105726872246696
$ws.Range("J16").Value = 0.3259911894273128
$ws.Range("K16").Value = 0.1233480176211454
$ws.Range("M16").Value = 0.03524229074889868
$ws.Range("O16").Value = 0.09251101321585903
$ws.Range("S16").Value = 0.1277533039647577
$ws.Range("F17").Value = 0.0233160621761658
$ws.Range("H17").Value = 0.2020725388601036
$ws.Range("I17").Value = 0.08808290155440414
$ws.Range("J17").Value = 0.4300518134715026
$ws.Range("K17").Value = 0.08549222797927461
$ws.Range("M17").Value = 0.02590673575129534
$ws.Range("O17").Value = 0.05958549222797927
$ws.Range("S17").Value = 0.08549222797927461
$ws.Range("F18").Value = 0.03061224489795918
$ws.Range("H18").Value = 0.1581632653061225
$ws.Range("I18").Value = 0.09183673469387756
$ws.Range("J18").Value = 0.4744897959183674
$ws.Range("K18").Value = 0.08163265306122448
$ws.Range("M18").Value = 0.01530612244897959
$ws.Range("O18").Value = 0.07142857142857142
$ws.Range("S18").Value = 0.07653061224489796
$ws.Range("F19").Value = 0.03267973856209151
$ws.Range("H19").Value = 0.2069716775599129
$ws.Range("J19").Value = 0.3710965867828613
$ws.Range("K19").Value = 0.1096586782861293
$ws.Range("M19").Value = 0.0196078431372549
$ws.Range("O19").Value = 0.06245461147421932
$ws.Range("S19").Value = 0.09949164851125636
$ws.Range("I19").Value = 0.09803921568627451
